# Update the "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly generated output data (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of row -> new F value
$updates = @{
    "展览" = @{
        3  = 592
        6  = 1120
        7  = 1464
        10 = 765
        12 = 185
        18 = 285
        19 = 5210
        20 = 88
        22 = 1021
        23 = 42
        26 = 6071
        31 = 14785
        36 = 10820
        37 = 667
        38 = 4235
    }
    "全部类型" = @{
        3  = 592
        6  = 1120
        7  = 1464
        10 = 765
        12 = 185
        18 = 285
        20 = 5210
        21 = 88
        24 = 1021
        25 = 42
        29 = 6071
        34 = 14785
        39 = 10820
        40 = 667
        41 = 4235
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
